# Add a new "phase lock" parameter/unit row (Painting / V) to the
# parameterUnit config sheet, right after the last existing row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A91").Value = "Painting"
$ws.Range("B91").Value = "V"

# Match the author's final selection in the sheet.
$ws.Range("C87").Select() | Out-Null
